$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "SSM" row (row 10): hours for "Number of new or modified input files" (I)
# went from 2 to 8, and a new "Actual (hours)" style figure of 8 was
# recorded in L10 (GWT Actual).
$ws.Range("I10").Value = 8
$ws.Range("L10").Value = 8

# Leave the cursor where the author left it when they saved the workbook.
$ws.Range("N14").Select()
